$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new BOM row for 0402 LED Green
$ws.Range("A4").Value = "0402 LED Green"

# Update selection to match target workbook view
$ws.Range("E42").Select()
